$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7, columns B (TB), C (d2S), D (K), E (IP), G (sum)
# Column F (Win) is unchanged.

$data = @{
    2 = @{ B = 0.01293466051926884; C = 0.04071648406533734; D = 0.1494219747398047;  E = 10.19245300693656;  G = 10.39552612626097 }
    3 = @{ B = 3.286832544864788;   C = 1.655778082260271;   D = 3.537761648806719;   E = 0.4942365360607697; G = 8.974608811992548 }
    4 = @{ B = 1.455362044514542;   C = 1.655778082260271;   D = 0.1494219747398047;  E = 10.19245300693656;  G = 13.45301510845117 }
    5 = @{ B = 1.455362044514542;   C = 1.655778082260271;   D = 3.537761648806719;   E = 10.19245300693656;  G = 16.84135478251809 }
    6 = @{ B = 3.286832544864788;   C = 1.655778082260271;   D = 0.7527432677738641;  E = 10.19245300693656;  G = 15.88780690183548 }
    7 = @{ B = 3.286832544864788;   C = 1.655778082260271;   D = 3.537761648806719;   E = 0.4942365360607697; G = 8.974608811992548 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
